$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update IPC PO (C), DELTA (D), DELTA^2 (E) for data rows 2-51,
# simulating weights reset to 0 so predictions are 0:
#   C = 0
#   D = -B (target - prediction)
#   E = D^2 (squared error)
$ws.Cells.Item(2, 3).Value = 0
$ws.Cells.Item(2, 4).Value = -29.55
$ws.Cells.Item(2, 5).Value = 873.2025
$ws.Cells.Item(3, 3).Value = 0
$ws.Cells.Item(3, 4).Value = -29.75
$ws.Cells.Item(3, 5).Value = 885.0625
$ws.Cells.Item(4, 3).Value = 0
$ws.Cells.Item(4, 4).Value = -29.84
$ws.Cells.Item(4, 5).Value = 890.4256
$ws.Cells.Item(5, 3).Value = 0
$ws.Cells.Item(5, 4).Value = -29.81
$ws.Cells.Item(5, 5).Value = 888.6360999999999
$ws.Cells.Item(6, 3).Value = 0
$ws.Cells.Item(6, 4).Value = -29.92
$ws.Cells.Item(6, 5).Value = 895.2064000000001
$ws.Cells.Item(7, 3).Value = 0
$ws.Cells.Item(7, 4).Value = -29.98
$ws.Cells.Item(7, 5).Value = 898.8004000000001
$ws.Cells.Item(8, 3).Value = 0
$ws.Cells.Item(8, 4).Value = -30.04
$ws.Cells.Item(8, 5).Value = 902.4015999999999
$ws.Cells.Item(9, 3).Value = 0
$ws.Cells.Item(9, 4).Value = -30.21
$ws.Cells.Item(9, 5).Value = 912.6441000000001
$ws.Cells.Item(10, 3).Value = 0
$ws.Cells.Item(10, 4).Value = -30.22
$ws.Cells.Item(10, 5).Value = 913.2484
$ws.Cells.Item(11, 3).Value = 0
$ws.Cells.Item(11, 4).Value = -30.38
$ws.Cells.Item(11, 5).Value = 922.9444
$ws.Cells.Item(12, 3).Value = 0
$ws.Cells.Item(12, 4).Value = -30.44
$ws.Cells.Item(12, 5).Value = 926.5936
$ws.Cells.Item(13, 3).Value = 0
$ws.Cells.Item(13, 4).Value = -30.48
$ws.Cells.Item(13, 5).Value = 929.0304
$ws.Cells.Item(14, 3).Value = 0
$ws.Cells.Item(14, 4).Value = -30.69
$ws.Cells.Item(14, 5).Value = 941.8761000000001
$ws.Cells.Item(15, 3).Value = 0
$ws.Cells.Item(15, 4).Value = -30.75
$ws.Cells.Item(15, 5).Value = 945.5625
$ws.Cells.Item(16, 3).Value = 0
$ws.Cells.Item(16, 4).Value = -30.94
$ws.Cells.Item(16, 5).Value = 957.2836000000001
$ws.Cells.Item(17, 3).Value = 0
$ws.Cells.Item(17, 4).Value = -30.95
$ws.Cells.Item(17, 5).Value = 957.9024999999999
$ws.Cells.Item(18, 3).Value = 0
$ws.Cells.Item(18, 4).Value = -31.02
$ws.Cells.Item(18, 5).Value = 962.2404
$ws.Cells.Item(19, 3).Value = 0
$ws.Cells.Item(19, 4).Value = -31.12
$ws.Cells.Item(19, 5).Value = 968.4544000000001
$ws.Cells.Item(20, 3).Value = 0
$ws.Cells.Item(20, 4).Value = -31.28
$ws.Cells.Item(20, 5).Value = 978.4384000000001
$ws.Cells.Item(21, 3).Value = 0
$ws.Cells.Item(21, 4).Value = -31.38
$ws.Cells.Item(21, 5).Value = 984.7044
$ws.Cells.Item(22, 3).Value = 0
$ws.Cells.Item(22, 4).Value = -31.58
$ws.Cells.Item(22, 5).Value = 997.2964
$ws.Cells.Item(23, 3).Value = 0
$ws.Cells.Item(23, 4).Value = -31.65
$ws.Cells.Item(23, 5).Value = 1001.7225
$ws.Cells.Item(24, 3).Value = 0
$ws.Cells.Item(24, 4).Value = -31.88
$ws.Cells.Item(24, 5).Value = 1016.3344
$ws.Cells.Item(25, 3).Value = 0
$ws.Cells.Item(25, 4).Value = -32.28
$ws.Cells.Item(25, 5).Value = 1041.9984
$ws.Cells.Item(26, 3).Value = 0
$ws.Cells.Item(26, 4).Value = -32.45
$ws.Cells.Item(26, 5).Value = 1053.0025
$ws.Cells.Item(27, 3).Value = 0
$ws.Cells.Item(27, 4).Value = -32.85
$ws.Cells.Item(27, 5).Value = 1079.1225
$ws.Cells.Item(28, 3).Value = 0
$ws.Cells.Item(28, 4).Value = -32.9
$ws.Cells.Item(28, 5).Value = 1082.41
$ws.Cells.Item(29, 3).Value = 0
$ws.Cells.Item(29, 4).Value = -33.1
$ws.Cells.Item(29, 5).Value = 1095.61
$ws.Cells.Item(30, 3).Value = 0
$ws.Cells.Item(30, 4).Value = -33.4
$ws.Cells.Item(30, 5).Value = 1115.56
$ws.Cells.Item(31, 3).Value = 0
$ws.Cells.Item(31, 4).Value = -33.7
$ws.Cells.Item(31, 5).Value = 1135.69
$ws.Cells.Item(32, 3).Value = 0
$ws.Cells.Item(32, 4).Value = -34.1
$ws.Cells.Item(32, 5).Value = 1162.81
$ws.Cells.Item(33, 3).Value = 0
$ws.Cells.Item(33, 4).Value = -34.4
$ws.Cells.Item(33, 5).Value = 1183.36
$ws.Cells.Item(34, 3).Value = 0
$ws.Cells.Item(34, 4).Value = -34.9
$ws.Cells.Item(34, 5).Value = 1218.01
$ws.Cells.Item(35, 3).Value = 0
$ws.Cells.Item(35, 4).Value = -35.3
$ws.Cells.Item(35, 5).Value = 1246.09
$ws.Cells.Item(36, 3).Value = 0
$ws.Cells.Item(36, 4).Value = -35.7
$ws.Cells.Item(36, 5).Value = 1274.49
$ws.Cells.Item(37, 3).Value = 0
$ws.Cells.Item(37, 4).Value = -36.3
$ws.Cells.Item(37, 5).Value = 1317.69
$ws.Cells.Item(38, 3).Value = 0
$ws.Cells.Item(38, 4).Value = -36.8
$ws.Cells.Item(38, 5).Value = 1354.24
$ws.Cells.Item(39, 3).Value = 0
$ws.Cells.Item(39, 4).Value = -37.3
$ws.Cells.Item(39, 5).Value = 1391.29
$ws.Cells.Item(40, 3).Value = 0
$ws.Cells.Item(40, 4).Value = -37.9
$ws.Cells.Item(40, 5).Value = 1436.41
$ws.Cells.Item(41, 3).Value = 0
$ws.Cells.Item(41, 4).Value = -38.5
$ws.Cells.Item(41, 5).Value = 1482.25
$ws.Cells.Item(42, 3).Value = 0
$ws.Cells.Item(42, 4).Value = -38.9
$ws.Cells.Item(42, 5).Value = 1513.21
$ws.Cells.Item(43, 3).Value = 0
$ws.Cells.Item(43, 4).Value = -39.4
$ws.Cells.Item(43, 5).Value = 1552.36
$ws.Cells.Item(44, 3).Value = 0
$ws.Cells.Item(44, 4).Value = -39.9
$ws.Cells.Item(44, 5).Value = 1592.01
$ws.Cells.Item(45, 3).Value = 0
$ws.Cells.Item(45, 4).Value = -40.1
$ws.Cells.Item(45, 5).Value = 1608.01
$ws.Cells.Item(46, 3).Value = 0
$ws.Cells.Item(46, 4).Value = -40.6
$ws.Cells.Item(46, 5).Value = 1648.36
$ws.Cells.Item(47, 3).Value = 0
$ws.Cells.Item(47, 4).Value = -40.9
$ws.Cells.Item(47, 5).Value = 1672.81
$ws.Cells.Item(48, 3).Value = 0
$ws.Cells.Item(48, 4).Value = -41.2
$ws.Cells.Item(48, 5).Value = 1697.44
$ws.Cells.Item(49, 3).Value = 0
$ws.Cells.Item(49, 4).Value = -41.5
$ws.Cells.Item(49, 5).Value = 1722.25
$ws.Cells.Item(50, 3).Value = 0
$ws.Cells.Item(50, 4).Value = -41.8
$ws.Cells.Item(50, 5).Value = 1747.24
$ws.Cells.Item(51, 3).Value = 0
$ws.Cells.Item(51, 4).Value = -42.2
$ws.Cells.Item(51, 5).Value = 1780.84

# TOTAL row: C52 sums the DELTA column, E52 sums DELTA^2
$ws.Cells.Item(52, 3).Value = -1702.24
$ws.Cells.Item(52, 5).Value = 58754.575

# MSE row: mean of DELTA^2 over the 50 data rows
$ws.Cells.Item(53, 5).Value = 1175.0915
